$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Sheet "Cutting Speed" (sheet1): add a new data row (row 4)
# ---------------------------------------------------------------------------
$wsCutting = $wb.Worksheets.Item(1)
$wsCutting.Range("B4").Value = 300
$wsCutting.Range("C4").Value = 80
$wsCutting.Range("D4").Value = 6
$wsCutting.Range("E4").Value = 0.12

# ---------------------------------------------------------------------------
# 2. Sheet "Material Removal Rate" (sheet2): update existing row 3, add row 4
# ---------------------------------------------------------------------------
$wsMrr = $wb.Worksheets.Item(2)
$wsMrr.Range("B3").Value = 5
$wsMrr.Range("C3").Value = 60
$wsMrr.Range("D3").Value = 1000
$wsMrr.Range("E3").Value = 300
$wsMrr.Range("F3").Value = "cm³/min"

$wsMrr.Range("B4").Value = 10
$wsMrr.Range("C4").Value = 60
$wsMrr.Range("D4").Value = 1000
$wsMrr.Range("E4").Value = 600
$wsMrr.Range("F4").Value = "cm³/min"

# ---------------------------------------------------------------------------
# 3. Sheet "Helix Angle" (sheet3): update existing row 3, add rows 4-5
# ---------------------------------------------------------------------------
$wsHelix = $wb.Worksheets.Item(3)
$wsHelix.Range("B3").Value = 5
$wsHelix.Range("C3").Value = 6
$wsHelix.Range("D3").Value = 0.2
$wsHelix.Range("E3").Value = 3.65
$wsHelix.Range("F3").Value = "°"

$wsHelix.Range("B4").Value = 5
$wsHelix.Range("C4").Value = 6
$wsHelix.Range("D4").Value = 0.1
$wsHelix.Range("E4").Value = 1.82
$wsHelix.Range("F4").Value = "°"

$wsHelix.Range("B5").Value = 5
$wsHelix.Range("C5").Value = 6
$wsHelix.Range("D5").Value = 0.06
$wsHelix.Range("E5").Value = 1.0900000000000001
$wsHelix.Range("F5").Value = "°"

# ---------------------------------------------------------------------------
# 4. New sheet "Ramp Angle" appended after "Helix Angle"
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsRamp = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$wsRamp.Name = "Ramp Angle"

$wsRamp.Range("B2").Value = "Toolpath Length"
$wsRamp.Range("C2").Value = "Step/Pitch"
$wsRamp.Range("D2").Value = "Angle of Decent"

$wsRamp.Range("B3").Value = 100
$wsRamp.Range("C3").Value = 1
$wsRamp.Range("D3").Value = 0.56999999999999995
$wsRamp.Range("E3").Value = "°"

$wsRamp.Range("B4").Value = 100
$wsRamp.Range("C4").Value = 3
$wsRamp.Range("D4").Value = 1.72
$wsRamp.Range("E4").Value = "°"

$wsRamp.Range("B5").Value = 150
$wsRamp.Range("C5").Value = 3
$wsRamp.Range("D5").Value = 1.1499999999999999
$wsRamp.Range("E5").Value = "°"

$wsRamp.Range("B6").Value = 200
$wsRamp.Range("C6").Value = 3
$wsRamp.Range("D6").Value = 0.86
$wsRamp.Range("E6").Value = "°"

# ---------------------------------------------------------------------------
# 5. New sheet "Surface Roughness" appended after "Ramp Angle" (becomes active)
# ---------------------------------------------------------------------------
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsRough = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet2)
$wsRough.Name = "Surface Roughness"

$wsRough.Range("B2").Value = "Feed per Turn"
$wsRough.Range("C2").Value = "Nose Radius"
$wsRough.Range("D2").Value = "Ra"

$wsRough.Range("B3").Value = 0.5
$wsRough.Range("C3").Value = 1
$wsRough.Range("D3").Value = 10.42

$wsRough.Range("B4").Value = 0.2
$wsRough.Range("C4").Value = 1
$wsRough.Range("D4").Value = 1.67

$wsRough.Range("B5").Value = 0.3
$wsRough.Range("C5").Value = 1
$wsRough.Range("D5").Value = 3.75

$wsRough.Range("B6").Value = 0.2
$wsRough.Range("C6").Value = 1
$wsRough.Range("D6").Value = 1.67

$wsRough.Range("B7").Value = 0.2
$wsRough.Range("C7").Value = 0.2
$wsRough.Range("D7").Value = 8.33

$wsRough.Range("B8").Value = 0.05
$wsRough.Range("C8").Value = 0.2
$wsRough.Range("D8").Value = 0.52

$wsRough.Range("B9").Value = 0.08
$wsRough.Range("C9").Value = 0.2
$wsRough.Range("D9").Value = 1.33

# ---------------------------------------------------------------------------
# 6. Clear the stale selection left on "Helix Angle" and make "Surface
#    Roughness" the active / selected sheet & cell, matching activeTab=4.
# ---------------------------------------------------------------------------
[void]$wsHelix.Range("A1").Select()
[void]$wsRough.Range("A1").Select()
